$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells with same style as existing header row (AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill team record (Wins/Losses/Ties) for every data row
for ($i = 2; $i -le 51; $i++) {
    $ws.Cells.Item($i, 30).Value = 92
    $ws.Cells.Item($i, 31).Value = 70
    $ws.Cells.Item($i, 32).Value = 0
}
